$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateSTP")

# Row 2-3: rename blank-field test data suffix from 04011/04012 to 07011
$ws.Range("C2").Value = "Blank Lead Contact07011"
$ws.Range("D2").Value = "BlankDescription07011"
$ws.Range("F2").Value = "FullName107011"
$ws.Range("G2").Value = "FullName207011"

$ws.Range("C3").Value = "LeadContBlank07011"
$ws.Range("D3").Value = "BlankDesc07011"
$ws.Range("E3").Value = "fullNameBlank07011"
$ws.Range("F3").Value = "Short107011"
$ws.Range("G3").Value = "Short207011"

# Rows 10-25, column G: rename suffix from 11 to 12
$ws.Range("G10").Value = "asset12"
$ws.Range("G11").Value = "technique12"
$ws.Range("G12").Value = "material12"
$ws.Range("G13").Value = "application12"
$ws.Range("G14").Value = "project12"
$ws.Range("G15").Value = "service12"
$ws.Range("G16").Value = "example12"
$ws.Range("G17").Value = "tech12"
$ws.Range("G18").Value = "Database12"
$ws.Range("G19").Value = "trend12"
$ws.Range("G20").Value = "Keyword12"
$ws.Range("G21").Value = "internal publi12"
$ws.Range("G22").Value = "external publi12"
$ws.Range("G23").Value = "collab12"
$ws.Range("G24").Value = "patent12"
$ws.Range("G25").Value = "compet12"

# Update the active selection to match the recorded state
$ws.Range("G10").Select()
